$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("calculations1")

# Lookup table used by VLOOKUP/HLOOKUP examples
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 4

$ws.Range("I7").Value = "a"
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 2

$ws.Range("I8").Value = "b"
$ws.Range("J8").Value = 3
$ws.Range("K8").Value = 4
$ws.Range("L8").Value = 5

$ws.Range("I9").Value = "c"
$ws.Range("J9").Value = 6
$ws.Range("K9").Value = 7
$ws.Range("L9").Value = 8

$ws.Range("I10").Value = "d"
$ws.Range("J10").Value = 9
$ws.Range("K10").Value = 10
$ws.Range("L10").Value = 11

# VLOOKUP / HLOOKUP examples
$ws.Range("B17").Value = "u"
$ws.Range("B18").Value = "t"

$ws.Range("C17").Formula = '=VLOOKUP("a",$I$6:$L$10,2,FALSE)'
$ws.Range("D17").Value = "SVERWEIS"

$ws.Range("C18").Formula = '=HLOOKUP(3,$I$6:$L$10,2,FALSE)'
$ws.Range("D18").Value = "WVERWEIS"

# Update selection to reflect where the user ended up working
[void]$ws.Range("D19").Select()
